$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 ("2021年") with the year's data, mirroring the layout of
# the existing rows (2012年..2020年) already present in the sheet.

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 20328.46
$ws.Range("C11").Value = 5034.95
$ws.Range("D11").Value = 868.73
# E11 intentionally left blank - no data reported for this period
$ws.Range("F11").Value = 17706.32
$ws.Range("G11").Value = 25949.22
$ws.Range("H11").Value = 3734.18
$ws.Range("I11").Value = 12169.42
$ws.Range("J11").Value = 3565.64
$ws.Range("K11").Value = 4153.8
$ws.Range("L11").Value = 2529.16
$ws.Range("M11").Value = 279.6
$ws.Range("N11").Value = 5514.92
$ws.Range("O11").Value = 13329.1
$ws.Range("P11").Value = 1339.54
$ws.Range("Q11").Value = 4200.29
$ws.Range("R11").Value = 13716.16
$ws.Range("S11").Value = 1012.77
$ws.Range("T11").Value = 16592.37
$ws.Range("U11").Value = 50.67
$ws.Range("V11").Value = 7389.54
$ws.Range("W11").Value = 1786.69
$ws.Range("X11").Value = 12159.33
$ws.Range("Y11").Value = 37656.45
$ws.Range("Z11").Value = 4104.81
$ws.Range("AA11").Value = 9977.4
$ws.Range("AB11").Value = 46.13
$ws.Range("AC11").Value = 409303.12
$ws.Range("AD11").Value = 13231.46
$ws.Range("AE11").Value = 5886.48
$ws.Range("AF11").Value = 30978.44
$ws.Range("AG11").Value = 22141.57
$ws.Range("AH11").Value = 5012.57
$ws.Range("AI11").Value = 4970.93
$ws.Range("AJ11").Value = 370.89
$ws.Range("AK11").Value = 22612.16
$ws.Range("AL11").Value = 4954.86
$ws.Range("AM11").Value = 36955.89
$ws.Range("AN11").Value = 2355.57
$ws.Range("AO11").Value = 7073.64
$ws.Range("AP11").Value = 24950.17
$ws.Range("AQ11").Value = 2612.05

# Copy the year-label formatting (bold, centered, bordered) from the cell
# above (A10) onto the new year label cell (A11), matching the style used
# by every other "year" cell in column A.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

# Touch E11 (without giving it a value) so a blank cell placeholder is
# emitted for it instead of it being omitted altogether - matching the
# source data which records an explicit (empty) cell for that column in
# this row.
$ws.Range("E11").Borders.LineStyle = 0

$excel.CutCopyMode = 0
